$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, pushing all existing rows (43..146) down
# by one (43->44 ... 146->147), matching the weekly data-refresh pattern
# seen in the diff (a new week's record is prepended, oldest record falls
# off the bottom into a brand-new row 147).
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with this week's record.
$ws.Cells.Item(43,1).Value  = 11
$ws.Cells.Item(43,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(43,3).Value  = "Bíobío"
$ws.Cells.Item(43,4).Value  = 45238
$ws.Cells.Item(43,5).Value  = 8
$ws.Cells.Item(43,6).Value  = 100112012
$ws.Cells.Item(43,7).Value  = "Espinaca"
$ws.Cells.Item(43,8).Value  = "Sin especificar"
$ws.Cells.Item(43,9).Value  = "Primera"
$ws.Cells.Item(43,10).Value = 50
$ws.Cells.Item(43,11).Value = 13000
$ws.Cells.Item(43,12).Value = 13000
$ws.Cells.Item(43,13).Value = 13000
$ws.Cells.Item(43,14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(43,15).Value = "Región Metropolitana"
$ws.Cells.Item(43,16).Value = 1300
$ws.Cells.Item(43,17).Value = 10
$ws.Cells.Item(43,18).Value = "Hortaliza"
